# Auto-generated Excel COM-interop edit script
# Applies updated crypto price/volume data per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value2 = "26.026.01"
$ws.Cells.Item(2, 5).Value2 = "  +0.45%  "
$ws.Cells.Item(3, 4).Value2 = "1.641.61"
$ws.Cells.Item(3, 5).Value2 = "  -0.45%  "
$ws.Cells.Item(4, 4).Value2 = "'1.001"
$ws.Cells.Item(4, 5).Value2 = "  -0.58%  "
$ws.Cells.Item(5, 4).Value2 = "'215.99"
$ws.Cells.Item(5, 5).Value2 = "  +0.15%  "
$ws.Cells.Item(6, 4).Value2 = "'0.5163"
$ws.Cells.Item(6, 5).Value2 = "  +1.44%  "
$ws.Cells.Item(7, 4).Value2 = "'1.002"
$ws.Cells.Item(7, 5).Value2 = "  -0.49%  "
$ws.Cells.Item(8, 5).Value2 = "  +0.45%  "
$ws.Cells.Item(9, 4).Value2 = "'0.06388"
$ws.Cells.Item(9, 5).Value2 = "  -0.64%  "
$ws.Cells.Item(10, 5).Value2 = "  +0.79%  "
$ws.Cells.Item(11, 4).Value2 = "'0.07780"
$ws.Cells.Item(11, 5).Value2 = "  -0.05%  "
$ws.Cells.Item(12, 5).Value2 = "  -0.53%  "
$ws.Cells.Item(13, 4).Value2 = "1.642.92"
$ws.Cells.Item(13, 5).Value2 = "  -0.03%  "
$ws.Cells.Item(14, 4).Value2 = "'0.5484"
$ws.Cells.Item(14, 5).Value2 = "  +0.09%  "
$ws.Cells.Item(15, 4).Value2 = "0.0{0}7800" -f [char]0x2085
$ws.Cells.Item(15, 5).Value2 = "  -1.41%  "
$ws.Cells.Item(16, 4).Value2 = "'64.63"
$ws.Cells.Item(16, 5).Value2 = "  -0.97%  "
$ws.Cells.Item(17, 4).Value2 = "26.036.26"
$ws.Cells.Item(17, 5).Value2 = "  +0.07%  "
$ws.Cells.Item(18, 5).Value2 = "  -0.48%  "
$ws.Cells.Item(19, 4).Value2 = "'199.90"
$ws.Cells.Item(19, 5).Value2 = "  +1.12%  "
$ws.Cells.Item(20, 4).Value2 = "'4.476"
$ws.Cells.Item(20, 5).Value2 = "  +0.83%  "
$ws.Cells.Item(21, 4).Value2 = "'10.00"
$ws.Cells.Item(21, 5).Value2 = "  -0.47%  "
$ws.Cells.Item(22, 4).Value2 = "'6.116"
$ws.Cells.Item(22, 5).Value2 = "  +0.76%  "
$ws.Cells.Item(23, 4).Value2 = "'1.003"
$ws.Cells.Item(23, 5).Value2 = "  -0.69%  "
$ws.Cells.Item(24, 4).Value2 = "'1.903"
$ws.Cells.Item(24, 5).Value2 = "  +2.42%  "
$ws.Cells.Item(25, 4).Value2 = "'142.43"
$ws.Cells.Item(25, 5).Value2 = "  +0.66%  "
$ws.Cells.Item(26, 5).Value2 = "  +7.41%  "
$ws.Cells.Item(27, 4).Value2 = "'6.889"
$ws.Cells.Item(27, 5).Value2 = "  -0.40%  "
$ws.Cells.Item(28, 4).Value2 = "'15.68"
$ws.Cells.Item(28, 5).Value2 = "  -0.48%  "
$ws.Cells.Item(29, 4).Value2 = "'1.246"
$ws.Cells.Item(29, 5).Value2 = "  +0.14%  "
$ws.Cells.Item(30, 4).Value2 = "'0.04870"
$ws.Cells.Item(30, 5).Value2 = "  -3.34%  "
$ws.Cells.Item(31, 4).Value2 = "'3.315"
$ws.Cells.Item(31, 5).Value2 = "  +1.05%  "
$ws.Cells.Item(32, 4).Value2 = "'3.244"
$ws.Cells.Item(32, 5).Value2 = "  +1.05%  "
$ws.Cells.Item(33, 4).Value2 = "'1.548"
$ws.Cells.Item(33, 5).Value2 = "  +0.06%  "
$ws.Cells.Item(34, 4).Value2 = "'2.382"
$ws.Cells.Item(34, 5).Value2 = "  +0.38%  "
$ws.Cells.Item(35, 4).Value2 = "'0.9212"
$ws.Cells.Item(35, 5).Value2 = "  +2.94%  "
$ws.Cells.Item(36, 4).Value2 = "'0.5605"
$ws.Cells.Item(37, 4).Value2 = "'2.571"
$ws.Cells.Item(37, 5).Value2 = "  -1.09%  "
$ws.Cells.Item(38, 4).Value2 = "1.120.21"
$ws.Cells.Item(38, 5).Value2 = "  -1.31%  "
$ws.Cells.Item(39, 4).Value2 = "'0.01579"
$ws.Cells.Item(39, 5).Value2 = "  +0.77%  "
$ws.Cells.Item(40, 4).Value2 = "'1.002"
$ws.Cells.Item(40, 5).Value2 = "  -0.73%  "
$ws.Cells.Item(41, 4).Value2 = "'2.532"
$ws.Cells.Item(41, 5).Value2 = "  -1.19%  "
$ws.Cells.Item(42, 4).Value2 = "'5.590"
$ws.Cells.Item(42, 5).Value2 = "  -1.52%  "
$ws.Cells.Item(43, 4).Value2 = "'0.8104"
$ws.Cells.Item(43, 5).Value2 = "  -0.81%  "
$ws.Cells.Item(44, 4).Value2 = "'99.85"
$ws.Cells.Item(44, 5).Value2 = "  -0.03%  "
$ws.Cells.Item(45, 5).Value2 = "  -1.30%  "
$ws.Cells.Item(46, 4).Value2 = "1.785.81"
$ws.Cells.Item(46, 5).Value2 = "  +0.00%  "
$ws.Cells.Item(47, 4).Value2 = "'0.4536"
$ws.Cells.Item(47, 5).Value2 = "  -0.08%  "
$ws.Cells.Item(48, 4).Value2 = "'55.46"
$ws.Cells.Item(48, 5).Value2 = "  +0.13%  "
$ws.Cells.Item(49, 5).Value2 = "  +0.01%  "
$ws.Cells.Item(50, 4).Value2 = "'0.05223"
$ws.Cells.Item(50, 5).Value2 = "  +2.49%  "
$ws.Cells.Item(51, 2).Value2 = "Algorand"
$ws.Cells.Item(51, 3).Value2 = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(51, 4).Value2 = "'0.09597"
$ws.Cells.Item(51, 5).Value2 = "  +0.20%  "
